$d = $word.ActiveDocument
$CR = [char]13

# Helper: build a minimal WordprocessingML package snippet containing a single
# paragraph with the given runs, for use with Range.InsertXML (classic Word
# COM "paste structured XML" behaviour). $runsXml is the already-built
# <w:r>...</w:r> markup for the paragraph.
function Get-ParaPackageXml([string]$runsXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# NOTE: paragraph objects/indices shift whenever the document is mutated, so
# every lookup below re-scans $d.Paragraphs by index right before use instead
# of caching stale Paragraph/Range references across edits.
function Find-ParagraphIndexByText($d, [string]$text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $t = $d.Paragraphs($i).Range.Text.TrimEnd($CR)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

function Find-ParagraphIndexStartingWith($d, [string]$prefix) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $t = $d.Paragraphs($i).Range.Text.TrimEnd($CR)
        if ($t.StartsWith($prefix)) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1) Delete the paragraph "the code generation. For a fully self contained
#    file" entirely (its text moved earlier in the list, see step 4 below).
# ---------------------------------------------------------------------------
$idx = Find-ParagraphIndexByText $d "the code generation. For a fully self contained file"
$d.Paragraphs($idx).Range.Delete()

# ---------------------------------------------------------------------------
# 2) Before the "Cleanup:" bullet (i.e. right after "RRN: ..."), insert a
#    brand-new ilvl-0 bullet "Back to schema loader:". Inserting *before*
#    "Cleanup:" means the new paragraph mark/run picks up clean, unbolded
#    formatting from the "Cleanup:" run rather than the bold "RRN" run.
# ---------------------------------------------------------------------------
$idx = Find-ParagraphIndexStartingWith $d "Cleanup:"
$d.Paragraphs($idx).Range.InsertParagraphBefore()
$newPara = $d.Paragraphs($idx)
$newPara.Range.ListFormat.ListLevelNumber = 1
$newParaBody = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$newParaBody.Text = "Back to schema loader:"

# ---------------------------------------------------------------------------
# 3) Replace the text of the first bullet ("Load the code from DB. ...")
#    with three separate runs.
# ---------------------------------------------------------------------------
$idx = Find-ParagraphIndexByText $d "Load the code from DB. doing all the queries of PG from the .net"
$p = $d.Paragraphs($idx)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = Get-ParaPackageXml('<w:r><w:t xml:space="preserve">the code generation. </w:t></w:r><w:r><w:t xml:space="preserve">Yey! </w:t></w:r><w:r><w:t>For a fully self contained file</w:t></w:r>')
$r.InsertXML($xml)

# ---------------------------------------------------------------------------
# 4) Replace the "RN: All table stuff..." bullet with a single plain run.
# ---------------------------------------------------------------------------
$idx = Find-ParagraphIndexStartingWith $d "RN:"
$p = $d.Paragraphs($idx)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = Get-ParaPackageXml('<w:r><w:t>work closely with old code. For starters, generate header. Take file name from command line, have a default file name if not</w:t></w:r>')
$r.InsertXML($xml)

# ---------------------------------------------------------------------------
# 5) Replace the "RRN: the _process functions..." bullet with a single plain
#    run.
# ---------------------------------------------------------------------------
$idx = Find-ParagraphIndexStartingWith $d "RRN:"
$p = $d.Paragraphs($idx)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = Get-ParaPackageXml('<w:r><w:t>go proc by proc on .net, just mimick that</w:t></w:r>')
$r.InsertXML($xml)

Write-Output "done"
